$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (model "m") - update formula text and recomputed stats
$ws.Range("C2").Value = "mean_accel ~ habitat_type * season * day_night + (1 | animal_id)"
$ws.Range("F2").Value = 0.107915285263934
$ws.Range("G2").Value = 97152.1299688854
$ws.Range("H2").Value = -194142.259937771
$ws.Range("K2").Value = -193389.203476531
$ws.Range("L2").Value = 80498

# Row 3
$ws.Range("I3").Value = 946.837889253657
$ws.Range("J3").Value = [double]"2.49324364133616e-206"

# Row 4
$ws.Range("I4").Value = 1686.34105142113

# Row 5
$ws.Range("F5").Value = 0.10835138099944
$ws.Range("G5").Value = 96220.8707511627
$ws.Range("H5").Value = -192397.741502325
$ws.Range("I5").Value = 1744.51843544553
$ws.Range("K5").Value = -192193.207648655

# Row 6
$ws.Range("F6").Value = 0.107922764853945
$ws.Range("G6").Value = 96025.4999556982
$ws.Range("H6").Value = -192038.999911396
$ws.Range("I6").Value = 2103.26002637448
$ws.Range("K6").Value = -191983.217951305

# Row 7
$ws.Range("F7").Value = 0.108436433145455
$ws.Range("G7").Value = 95845.748974635
$ws.Range("I7").Value = 2462.7619885008
$ws.Range("K7").Value = -191623.715989178

# Row 8
$ws.Range("F8").Value = 0.108084000859798
$ws.Range("G8").Value = 95419.8018318772
$ws.Range("H8").Value = -190825.603663754
$ws.Range("I8").Value = 3316.65627401657
$ws.Range("K8").Value = -190760.524710314
